$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared string "Mao-de-obra " -> "MO" (cell D1 on sheet1)
$ws.Range("D1").Value = "MO"

# Update selection to G2
$ws.Range("G2").Select()

# Remove custom row height on row 1 (auto-fit)
$ws.Rows("1").AutoFit()
